$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 372.75
$ws.Range("I42").Value = 371.55554
$ws.Range("J42").Value = 376.33334
$ws.Range("K42").Value = 1114.66662
$ws.Range("L42").Value = 1129.00002
$ws.Range("M42").Value = -884.66662
$ws.Range("N42").Value = -1589.00002
$ws.Range("H43").Value = 3417.2
$ws.Range("I43").Value = 1480.2
$ws.Range("K43").Value = 1480.2
$ws.Range("M43").Value = -1411.2
$ws.Range("H58").Value = 1908.2106
$ws.Range("I58").Value = 250.54546
$ws.Range("J58").Value = 4187.5
$ws.Range("K58").Value = 751.6363799999999
$ws.Range("L58").Value = 12562.5
$ws.Range("M58").Value = -601.6363799999999
$ws.Range("N58").Value = -12862.5
$ws.Range("H82").Value = 3481.6667
$ws.Range("I82").Value = 222.5
$ws.Range("K82").Value = 667.5
$ws.Range("M82").Value = -261.5
$ws.Range("H85").Value = 3481.6667
$ws.Range("I85").Value = 222.5
$ws.Range("K85").Value = 667.5
$ws.Range("M85").Value = 736.5
$ws.Range("H99").Value = 296
$ws.Range("J99").Value = 499
$ws.Range("L99").Value = 1497
$ws.Range("N99").Value = -4493
$ws.Range("H101").Value = 235
$ws.Range("I101").Value = 241.875
$ws.Range("K101").Value = 725.625
$ws.Range("M101").Value = 896.375
$ws.Range("H118").Value = 964.06665
$ws.Range("I118").Value = 634.7692
$ws.Range("J118").Value = 3104.5
$ws.Range("K118").Value = 1904.3076
$ws.Range("L118").Value = 9313.5
$ws.Range("M118").Value = -247.3075999999999
$ws.Range("N118").Value = -12627.5
$ws.Range("H129").Value = 764.6667
$ws.Range("I129").Value = 610.25
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 1830.75
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 3169.25
$ws.Range("N129").Value = -16000
$ws.Range("H137").Value = 31668.9
$ws.Range("I137").Value = 43713.43
$ws.Range("K137").Value = 131140.29
$ws.Range("M137").Value = -128590.29
$ws.Range("H138").Value = 16860.193
$ws.Range("J138").Value = 44750.082
$ws.Range("L138").Value = 134250.246
$ws.Range("N138").Value = -144530.246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23646.223
$ws.Range("I32").Value = 24687.93
$ws.Range("K32").Value = 24687.93
$ws.Range("M32").Value = -24400.93
$ws.Range("H122").Value = 3687.3704
$ws.Range("I122").Value = 3637.95
$ws.Range("K122").Value = 10913.85
$ws.Range("M122").Value = -8463.849999999999
$ws.Range("H132").Value = 1511.1072
$ws.Range("I132").Value = 1232.44
$ws.Range("J132").Value = 3833.3333
$ws.Range("K132").Value = 3697.32
$ws.Range("L132").Value = 11499.9999
$ws.Range("M132").Value = -1167.32
$ws.Range("N132").Value = -16559.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 19801.176
$ws.Range("J20").Value = 1957.125
$ws.Range("L20").Value = 1957.125
$ws.Range("N20").Value = -2451.125
$ws.Range("H86").Value = 1753.6923
$ws.Range("I86").Value = 1649.8334
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 1649.8334
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -526.8334
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 1753.6923
$ws.Range("I89").Value = 1649.8334
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 8249.166999999999
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -2633.166999999999
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1069.7059
$ws.Range("I5").Value = 822
$ws.Range("J5").Value = 1874.75
$ws.Range("K5").Value = 2466
$ws.Range("L5").Value = 5624.25
$ws.Range("M5").Value = -2354
$ws.Range("N5").Value = -5848.25
$ws.Range("H135").Value = 1069.7059
$ws.Range("I135").Value = 822
$ws.Range("J135").Value = 1874.75
$ws.Range("K135").Value = 7398
$ws.Range("L135").Value = 16872.75
$ws.Range("M135").Value = -4863
$ws.Range("N135").Value = -21942.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8718.833000000001
$ws.Range("J70").Value = 7546.3335
$ws.Range("L70").Value = 7546.3335
$ws.Range("N70").Value = -8086.3335
$ws.Range("H73").Value = 8718.833000000001
$ws.Range("J73").Value = 7546.3335
$ws.Range("L73").Value = 7546.3335
$ws.Range("M73").Value = -9191.5
$ws.Range("N73").Value = -9418.333500000001
$ws.Range("H80").Value = 10298.667
$ws.Range("I80").Value = 4448
$ws.Range("J80").Value = 22000
$ws.Range("K80").Value = 4448
$ws.Range("L80").Value = 22000
$ws.Range("M80").Value = -3450
$ws.Range("N80").Value = -23996
$ws.Range("H83").Value = 10298.667
$ws.Range("I83").Value = 4448
$ws.Range("J83").Value = 22000
$ws.Range("K83").Value = 22240
$ws.Range("L83").Value = 110000
$ws.Range("M83").Value = -17248
$ws.Range("N83").Value = -119984
$ws.Range("H102").Value = 3200.682
$ws.Range("I102").Value = 2245.0715
$ws.Range("K102").Value = 2245.0715
$ws.Range("M102").Value = -623.0715
$ws.Range("H122").Value = 3224.3928
$ws.Range("I122").Value = 2959.6191
$ws.Range("K122").Value = 8878.8573
$ws.Range("M122").Value = -6428.8573
$ws.Range("H132").Value = 2989.8333
$ws.Range("I132").Value = 2989.8333
$ws.Range("K132").Value = 8969.499899999999
$ws.Range("M132").Value = -6439.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 818.4737
$ws.Range("I16").Value = 877.41174
$ws.Range("K16").Value = 877.41174
$ws.Range("M16").Value = -707.41174
$ws.Range("H22").Value = 1449.25
$ws.Range("I22").Value = 1095.0667
$ws.Range("J22").Value = 2039.5555
$ws.Range("K22").Value = 1095.0667
$ws.Range("L22").Value = 2039.5555
$ws.Range("M22").Value = -800.0667000000001
$ws.Range("N22").Value = -2629.5555
$ws.Range("H27").Value = 1449.25
$ws.Range("I27").Value = 1095.0667
$ws.Range("J27").Value = 2039.5555
$ws.Range("K27").Value = 1095.0667
$ws.Range("L27").Value = 2039.5555
$ws.Range("M27").Value = -988.0667000000001
$ws.Range("N27").Value = -2253.5555
$ws.Range("H40").Value = 2417.5334
$ws.Range("I40").Value = 2289.4614
$ws.Range("J40").Value = 3250
$ws.Range("K40").Value = 2289.4614
$ws.Range("L40").Value = 3250
$ws.Range("M40").Value = -2153.4614
$ws.Range("N40").Value = -3522
$ws.Range("H55").Value = 1137.3077
$ws.Range("I55").Value = 299.66666
$ws.Range("J55").Value = 1855.2858
$ws.Range("K55").Value = 299.66666
$ws.Range("L55").Value = 1855.2858
$ws.Range("M55").Value = -126.66666
$ws.Range("N55").Value = -2201.2858
$ws.Range("H61").Value = 927.2353000000001
$ws.Range("I61").Value = 672.6875
$ws.Range("K61").Value = 672.6875
$ws.Range("M61").Value = -470.6875
$ws.Range("H68").Value = 5450
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 5900
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 5900
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -7398
$ws.Range("H71").Value = 5450
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 5900
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 29500
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -36988
$ws.Range("H100").Value = 2796.158
$ws.Range("I100").Value = 2661.5334
$ws.Range("J100").Value = 3301
$ws.Range("K100").Value = 2661.5334
$ws.Range("L100").Value = 3301
$ws.Range("M100").Value = -2120.5334
$ws.Range("N100").Value = -4383
$ws.Range("H113").Value = 927.2353000000001
$ws.Range("I113").Value = 672.6875
$ws.Range("K113").Value = 672.6875
$ws.Range("M113").Value = 1497.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 613.0526
$ws.Range("I100").Value = 611.05884
$ws.Range("K100").Value = 1222.11768
$ws.Range("M100").Value = -681.1176800000001
$ws.Range("H107").Value = 1279.25
$ws.Range("J107").Value = 1371.25
$ws.Range("L107").Value = 4113.75
$ws.Range("N107").Value = -7953.75
$ws.Range("H113").Value = 491.5
$ws.Range("J113").Value = 2221.5
$ws.Range("L113").Value = 6664.5
$ws.Range("N113").Value = -11004.5
$ws.Range("H122").Value = 49105.3
$ws.Range("I122").Value = 58147.08
$ws.Range("K122").Value = 174441.24
$ws.Range("M122").Value = -171991.24
$ws.Range("H132").Value = 1100
$ws.Range("I132").Value = 927.88464
$ws.Range("K132").Value = 2783.65392
$ws.Range("M132").Value = -253.6539199999997
